$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New product row (row 8): LB318 bandolera, 5 image variants, 5 colors ---
$ws.Range("A8").Value = "LB318"
$ws.Range("B8").Value = "Bandolera"
$ws.Range("C8").Value = "Bandolera"
$ws.Range("D8").Value = "Marroquineria"
$ws.Range("E8").Value = 695
$ws.Range("F8").Value = "25*12*19"
$ws.Range("G8").Value = "LB318-1"
$ws.Range("H8").Value = "LB318-2"
$ws.Range("I8").Value = "LB318-3"
$ws.Range("J8").Value = "LB318-4"
$ws.Range("K8").Value = "LB318-5"

# Match the formatting already used by the "Imagen" columns on the rows above
$ws.Range("G8:K8").Style = $ws.Range("G7").Style

# Array formula that builds the CDN URL for the first image variant, same
# pattern used by every row above (R column = "Imagen 1 URL")
$ws.Range("R8").FormulaArray = '=IF(INDEX($G8:$Q8, COLUMN(A7))="", "", "https://cdn.jsdelivr.net/gh/Ferabensrl/catalogo-mare@main/imagenes/" & INDEX($G8:$Q8, COLUMN(A7)))'
$ws.Range("R8").Style = $ws.Range("R7").Style

# Color availability flags: Negro, Rosado, Beige, Bordeaux, Rosa Viejo
$ws.Range("AC8").Value = "SI"
$ws.Range("AP8").Value = "SI"
$ws.Range("AS8").Value = "SI"
$ws.Range("AY8").Value = "SI"
$ws.Range("BA8").Value = "SI"

$excel.Calculate()

# Widen the "Imagen 1 URL" column (R) so the longer URLs fit
$ws.Columns.Item(18).ColumnWidth = 76.42

# Freeze column A and move the selection to the newly added row's "Rosado" cell
$ws.Range("B1").Select()
$win = $ws.Application.ActiveWindow
$win.FreezePanes = $true
$ws.Range("AP8").Select()
